$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WSM - Cont. Orchestration Sys.")
$ws.Range("E23").Font.FontStyle = "Regular"
